# Fruta / hortaliza, semanal
#
# Insert two new weekly price rows for "Plátano" (Pintón / Primera Pintón)
# dated 2021-09-22 (serial 44461) right before the existing row 320 block,
# pushing all the subsequent rows (previously 320:346) down to (322:348).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 320 (existing rows 320.. shift down to 322..)
$ws.Rows.Item(320).Insert()
$ws.Rows.Item(320).Insert()

# New row 320: "Pintón"
$ws.Range("A320").Value = 7
$ws.Range("B320").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C320").Value = "Ñuble"
$ws.Range("D320").Value = 44461
$ws.Range("E320").Value = 16
$ws.Range("F320").Value = "Fruta"
$ws.Range("G320").Value = 100108
$ws.Range("H320").Value = "Tropicales y subtropicales"
$ws.Range("I320").Value = 100108006
$ws.Range("J320").Value = "Plátano"
$ws.Range("K320").Value = "Sin especificar"
$ws.Range("L320").Value = "Pintón"
$ws.Range("M320").Value = 150
$ws.Range("N320").Value = 16000
$ws.Range("O320").Value = 16000
$ws.Range("P320").Value = 16000
$ws.Range("Q320").Value = "$/caja 20 kilos"
$ws.Range("R320").Value = "Ecuador"
$ws.Range("S320").Value = 800
$ws.Range("T320").Value = 20

# New row 321: "Primera Pintón"
$ws.Range("A321").Value = 7
$ws.Range("B321").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C321").Value = "Ñuble"
$ws.Range("D321").Value = 44461
$ws.Range("E321").Value = 16
$ws.Range("F321").Value = "Fruta"
$ws.Range("G321").Value = 100108
$ws.Range("H321").Value = "Tropicales y subtropicales"
$ws.Range("I321").Value = 100108006
$ws.Range("J321").Value = "Plátano"
$ws.Range("K321").Value = "Sin especificar"
$ws.Range("L321").Value = "Primera Pintón"
$ws.Range("M321").Value = 300
$ws.Range("N321").Value = 17000
$ws.Range("O321").Value = 18000
$ws.Range("P321").Value = 17500
$ws.Range("Q321").Value = "$/caja 20 kilos"
$ws.Range("R321").Value = "Ecuador"
$ws.Range("S321").Value = 875
$ws.Range("T321").Value = 20

Write-Host ("New dimension: " + $ws.UsedRange.Address())
